$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Direct text/string assignments (not numeric-ambiguous) ---
$ws.Range("D2").Value = '26.431.74'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.837.94'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("E7").Value = '  +0.95%  '
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '1.849.45'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '26.465.62'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = '2.082.75'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("E26").Value = '  -4.43%  '
$ws.Range("E27").Value = '  +2.10%  '
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  -2.59%  '
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E38").Value = '  +1.26%  '
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("E40").Value = '  -2.82%  '
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  -1.37%  '
$ws.Range("E43").Value = '  -3.05%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  +0.62%  '

# --- Numeric-looking values that must remain TEXT: use helper column Z ---
#     (apostrophe-prefixed helper cell copied via PasteSpecial values-only
#      avoids Excel auto-converting these to numbers, while keeping the
#      destination cell at default style, matching the source formatting.)
$ws.Range("Z2").Value = "'1.002"
$ws.Range("Z2").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("Z3").Value = "'260.04"
$ws.Range("Z3").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z4").Value = "'0.5259"
$ws.Range("Z4").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("Z5").Value = "'0.3200"
$ws.Range("Z5").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("Z6").Value = "'0.06780"
$ws.Range("Z6").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("Z7").Value = "'18.73"
$ws.Range("Z7").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z8").Value = "'0.7812"
$ws.Range("Z8").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("Z9").Value = "'0.07756"
$ws.Range("Z9").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("Z10").Value = "'87.52"
$ws.Range("Z10").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("Z11").Value = "'4.999"
$ws.Range("Z11").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("Z12").Value = "'1.002"
$ws.Range("Z12").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("Z13").Value = "'13.82"
$ws.Range("Z13").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("Z14").Value = "'1.002"
$ws.Range("Z14").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("Z15").Value = "'0.000007913"
$ws.Range("Z15").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("Z16").Value = "'4.608"
$ws.Range("Z16").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("Z17").Value = "'9.360"
$ws.Range("Z17").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("Z18").Value = "'5.954"
$ws.Range("Z18").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("Z19").Value = "'141.42"
$ws.Range("Z19").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("Z20").Value = "'2.184"
$ws.Range("Z20").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("Z21").Value = "'16.90"
$ws.Range("Z21").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("Z22").Value = "'111.60"
$ws.Range("Z22").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("Z23").Value = "'4.150"
$ws.Range("Z23").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("Z24").Value = "'4.050"
$ws.Range("Z24").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("Z25").Value = "'0.04862"
$ws.Range("Z25").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z26").Value = "'0.7241"
$ws.Range("Z26").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("Z27").Value = "'1.128"
$ws.Range("Z27").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("Z28").Value = "'2.861"
$ws.Range("Z28").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("Z29").Value = "'3.090"
$ws.Range("Z29").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("Z30").Value = "'2.231"
$ws.Range("Z30").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("Z31").Value = "'0.01767"
$ws.Range("Z31").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("Z32").Value = "'0.4757"
$ws.Range("Z32").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("Z33").Value = "'0.8902"
$ws.Range("Z33").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("Z34").Value = "'109.60"
$ws.Range("Z34").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("Z35").Value = "'5.908"
$ws.Range("Z35").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("Z36").Value = "'1.002"
$ws.Range("Z36").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("Z37").Value = "'7.645"
$ws.Range("Z37").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z38").Value = "'0.4125"
$ws.Range("Z38").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("Z39").Value = "'9.007"
$ws.Range("Z39").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z40").Value = "'0.05858"
$ws.Range("Z40").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("Z41").Value = "'0.1228"
$ws.Range("Z41").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("Z42").Value = "'34.84"
$ws.Range("Z42").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("Z43").Value = "'0.8900"
$ws.Range("Z43").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$ws.Range("Z2:Z43").Clear()
$excel.CutCopyMode = $false
